# This script updates a handful of numeric "view/favorite count" style
# values (column F) across the workbook's sheets, matching the
# gh-pages data refresh captured in the target diff.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (cell, new value) updates
$updates = @{
    "展览"   = @(
        @{ Cell = "F2";  Value = 27 },
        @{ Cell = "F7";  Value = 940 },
        @{ Cell = "F11"; Value = 1409 },
        @{ Cell = "F14"; Value = 2974 },
        @{ Cell = "F15"; Value = 370 },
        @{ Cell = "F18"; Value = 775 },
        @{ Cell = "F20"; Value = 1338 },
        @{ Cell = "F25"; Value = 3413 },
        @{ Cell = "F26"; Value = 667 }
    )
    "演出"   = @(
        @{ Cell = "F3";  Value = 45 },
        @{ Cell = "F10"; Value = 5 }
    )
    "全部类型" = @(
        @{ Cell = "F2";  Value = 27 },
        @{ Cell = "F8";  Value = 45 },
        @{ Cell = "F9";  Value = 45 },
        @{ Cell = "F17"; Value = 940 },
        @{ Cell = "F21"; Value = 1409 },
        @{ Cell = "F24"; Value = 2974 },
        @{ Cell = "F25"; Value = 370 },
        @{ Cell = "F28"; Value = 775 },
        @{ Cell = "F30"; Value = 1338 },
        @{ Cell = "F33"; Value = 5 },
        @{ Cell = "F37"; Value = 3413 },
        @{ Cell = "F38"; Value = 667 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates[$sheetName]) {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
